# MAJOR RELEASE: MACHINE LEARNING MODEL V1.0
# Swap the row data (columns B:K) between row 16 (Ron Harper Jr. (TW))
# and row 17 (Will Barton), leaving column A ("No.") untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($col = 2; $col -le 11; $col++) {
    $cellA = $ws.Cells.Item(16, $col)
    $cellB = $ws.Cells.Item(17, $col)

    $tmp = $cellA.Value()
    $valB = $cellB.Value()

    # Column I ("Exp") can hold text that looks numeric (e.g. "10" years
    # of experience, stored as a shared string, not a literal number).
    # Preserve the original per-cell style (so no extra formatting is
    # introduced) while forcing the re-assigned value to stay textual.
    if ($col -eq 9) {
        $styleA = $cellA.Style
        $styleB = $cellB.Style

        $cellA.NumberFormat = "@"
        $cellA.Value = $valB
        $cellA.Style = $styleA

        $cellB.NumberFormat = "@"
        $cellB.Value = $tmp
        $cellB.Style = $styleB
    } else {
        $cellA.Value = $valB
        $cellB.Value = $tmp
    }
}
